# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# The "K" column (column G) values were recalculated; update the new values in place.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 0
    3  = 0
    4  = 1
    5  = 1
    6  = 0
    7  = 2
    8  = 1
    9  = 0
    10 = 3
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
